$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3016
$ws.Range("C2").Value = 0.3016

$ws.Range("B3").Value = 0.2996
$ws.Range("C3").Value = 0.2996

$ws.Range("B4").Value = 0.2966
$ws.Range("C4").Value = 0.2966

$ws.Range("B5").Value = 0.2912
$ws.Range("C5").Value = 0.2912
